# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# Changes applied:
#  1. About sheet (sheet 1): update the "last updated" date in C1 from
#     2024-03-15 (45366) to 2024-03-28 (45379).
#  2. RAF-capacity sheet (sheet 4): raise the capacity-credit multiplier for
#     the two hydrogen technologies (rows 24-25, "hydrogen combustion
#     turbine" / "hydrogen combined cycle") from 0.3 to 1.
#  3. View-state changes captured in the saved workbook: RAF-capacity
#     becomes the active/selected sheet (instead of RAF-generation), its
#     zoom is set to 80%, its column A is widened to fit the longer labels,
#     and the last-used cell selection on that sheet moves to B25.

$wb = $excel.ActiveWorkbook

$wsAbout       = $wb.Worksheets.Item(1)   # About
$wsGeneration  = $wb.Worksheets.Item(2)   # RAF-generation
$wsDemandAlter = $wb.Worksheets.Item(3)   # RAF-demand-altering-techs
$wsCapacity    = $wb.Worksheets.Item(4)   # RAF-capacity

# 1. Bump the "Updated" date on the About sheet.
$wsAbout.Range("C1").Value = "2024-03-28"

# 2. Raise the hydrogen plant capacity-credit multipliers to 1.
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# 3. Widen column A on RAF-capacity so the hydrogen tech labels fit, and
#    make it the active sheet / selection the workbook opens on.
$wsCapacity.Columns.Item(1).ColumnWidth = 28.1

$wsCapacity.Activate()
$excel.ActiveWindow.Zoom = 80
$wsCapacity.Range("B25").Select()
